$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '28.611.52'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  -3.07%  '

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '1.850.76'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  -3.83%  '

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '1.002'
$cell.Style = "Normal"
$ws.Range("E4").Value = '  -0.97%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '336.20'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +3.03%  '

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '1.002'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  -0.93%  '

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.4665'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  -3.01%  '

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.3902'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  -3.60%  '

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '46.16'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  -2.36%  '

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '0.07916'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  -3.36%  '

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.9803'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  -2.71%  '

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '22.29'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  -6.25%  '

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '1.829.39'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  -5.09%  '

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '5.828'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  -4.20%  '

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '6.993'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  -4.20%  '

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '0.06910'
$cell.Style = "Normal"

$ws.Range("E17").Value = '  -1.06%  '

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '87.63'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  -4.24%  '

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '0.00001003'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  -3.27%  '

$ws.Range("E20").Value = '  -2.97%  '

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '1.002'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  -0.82%  '

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '28.617.69'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  -3.09%  '

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '5.393'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  -4.73%  '

$ws.Range("E24").Value = '  -5.77%  '

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '2.165'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  -0.69%  '

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '2.097.07'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  -3.17%  '

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '153.37'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  -1.61%  '

$ws.Range("E28").Value = '  -2.97%  '

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '6.062'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  -5.08%  '

$ws.Range("E30").Value = '  -2.80%  '

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '117.40'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  -2.55%  '

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '0.9715'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  -4.04%  '

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '0.09339'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  -2.52%  '

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '5.357'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  -4.23%  '

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '3.485'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  -2.09%  '

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '1.346'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  -2.81%  '

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '0.06159'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -2.97%  '

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '0.02201'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  -3.35%  '

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '1.171'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  -1.51%  '

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '10.15'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  -5.28%  '

$ws.Range("E43").Value = '  -2.63%  '

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '2.415'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  -2.30%  '

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '1.249'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  -2.11%  '

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '0.07102'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  -4.91%  '

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '1.905'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  -3.41%  '

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '113.77'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  -3.23%  '

# Row 40 full update
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '7.691'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  -2.50%  '

# Row 41 full update
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '0.5709'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  -3.94%  '

# Row 46 full update
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '11.82'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  -4.75%  '

# Row 47 full update
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '0.5385'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  -2.88%  '

# Row 51 full update
$ws.Range("B51").Value = 'PaxDollar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '1.002'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  -0.93%  '
